$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Code status" (column F) cells -------------------------------
# Rows moving from "IN PROGRESS" -> "DONE"
$doneSource = $ws.Range("F2")
foreach ($addr in @("F6", "F37", "F38")) {
    $doneSource.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($addr).Value = "DONE"
}

# Row 22 moving from "TODO" -> "DONE"
$doneSource.Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("F22").Value = "DONE"

# Rows moving from "TODO" -> "IN PROGRESS: EN" (new shared string added first)
$inProgressSource = $ws.Range("F5")
foreach ($addr in @("F62", "F63")) {
    $inProgressSource.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = "IN PROGRESS: EN"
}

# Rows moving from "TODO" -> "IN PROGRESS: HK" (new shared string added second)
foreach ($addr in @("F49", "F50")) {
    $inProgressSource.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = "IN PROGRESS: HK"
}

$excel.CutCopyMode = 0

# --- Update the active selection shown when the workbook is reopened ----
$ws.Range("F58").Select() | Out-Null
